# Update plots for each sample: refresh peak-detection results after
# re-running the caller on updated bin ranges for a couple of markers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# peak_table: a couple of bin boundaries / a height were recomputed
# ---------------------------------------------------------------------
$peak = $wb.Worksheets.Item("peak_table")

# S1 / CYP2D6_49 (row 4): w_height dropped from 1000 to 800
$peak.Range("N4").Value = 800

# S2 / CYP2D6_4 (row 12): w_max widened from 32 to 33
$peak.Range("G12").Value = 33

# S2 / CYP2D6_3 (row 13): m_min widened from 32 to 33
$peak.Range("H13").Value = 33

# ---------------------------------------------------------------------
# allele_table: per-allele detection results follow the new bins
# ---------------------------------------------------------------------
$allele = $wb.Worksheets.Item("allele_table")

# Row 6: S1 / CYP2D6_003 (CYP2D6_49) wildtype "T" - now detected
$allele.Range("K6").Value = 800
$allele.Range("M6").Value = $true
$allele.Range("N6").Value = 25
$allele.Range("O6").Value = 39.76
$allele.Range("P6").Value = 994
$allele.Range("Q6").Value = "ok"
$allele.Range("R6").Value = ""

# Row 22: S2 / CYP2D6_011 (CYP2D6_4) wildtype "G" - now detected
$allele.Range("J22").Value = 33
$allele.Range("M22").Value = $true
$allele.Range("N22").Value = 54
$allele.Range("O22").Value = 32.26
$allele.Range("P22").Value = 1082
$allele.Range("Q22").Value = "ok"
$allele.Range("R22").Value = ""

# Row 25: S2 / CYP2D6_012 (CYP2D6_3) mutant "G" - no longer detected
$allele.Range("I25").Value = 33
$allele.Range("M25").Value = $false
$allele.Range("N25").Value = ""
$allele.Range("O25").Value = ""
$allele.Range("P25").Value = ""
$allele.Range("Q25").Value = ""
$allele.Range("R25").Value = "Peak(s) could not be detected. Please check peak ranges if required!"

# ---------------------------------------------------------------------
# marker_table: genotype/phenotype calls follow allele detection
# ---------------------------------------------------------------------
$marker = $wb.Worksheets.Item("marker_table")

# Row 4: S1 / CYP2D6_003 (CYP2D6_49) now resolves to a TT wildtype call
$marker.Range("G4").Value = "TT"
$marker.Range("H4").Value = "wildtype"

# Row 12: S2 / CYP2D6_011 (CYP2D6_4) now resolves to a GG wildtype call
$marker.Range("G12").Value = "GG"
$marker.Range("H12").Value = "wildtype"

# Row 13: S2 / CYP2D6_012 (CYP2D6_3) mutant allele no longer detected -> wildtype
$marker.Range("G13").Value = "AA"
$marker.Range("H13").Value = "wildtype"

# ---------------------------------------------------------------------
# genotype_result: final diplotype call for the sample
# ---------------------------------------------------------------------
$genotype = $wb.Worksheets.Item("genotype_result")
$genotype.Range("B2").Value = "*1/*2"
